$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.527126666666668
$ws.Range("H2").Value = 19.58138
$ws.Range("I2").Value = 0.2130391554800433
$ws.Range("J2").Value = 0.2130391554800433
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.62053
$ws.Range("N2").Value = 22.86159
$ws.Range("O2").Value = 0.2214906134136664
$ws.Range("P2").Value = 0.2214906134136664
$ws.Range("Q2").Value = 49.74016457713333
$ws.Range("R2").Value = 447.6614811942001
$ws.Range("S2").Value = 0.04718617322840424
$ws.Range("T2").Value = 0.04718617322840424
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.527126666666668
$ws.Range("H3").Value = 19.58138
$ws.Range("I3").Value = 0.2130391554800433
$ws.Range("J3").Value = 0.2130391554800433
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.340016666666665
$ws.Range("N3").Value = 28.02005
$ws.Range("O3").Value = 0.2714674728390108
$ws.Range("P3").Value = 0.2714674728390108
$ws.Range("Q3").Value = 60.96347185211111
$ws.Range("R3").Value = 548.6712466690001
$ws.Range("S3").Value = 0.05783320115392448
$ws.Range("T3").Value = 0.05783320115392448
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.527126666666668
$ws.Range("H4").Value = 19.58138
$ws.Range("I4").Value = 0.2130391554800433
$ws.Range("J4").Value = 0.2130391554800433
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.970095666666666
$ws.Range("N4").Value = 23.910287
$ws.Range("O4").Value = 0.2316507353393536
$ws.Range("P4").Value = 0.2316507353393536
$ws.Range("Q4").Value = 52.02182396178445
$ws.Range("R4").Value = 468.1964156560601
$ws.Range("S4").Value = 0.04935067702302692
$ws.Range("T4").Value = 0.04935067702302692
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.527126666666668
$ws.Range("H5").Value = 19.58138
$ws.Range("I5").Value = 0.2130391554800433
$ws.Range("J5").Value = 0.2130391554800433
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.475014333333332
$ws.Range("N5").Value = 28.425043
$ws.Range("O5").Value = 0.2753911784079691
$ws.Range("P5").Value = 0.2753911784079691
$ws.Range("Q5").Value = 61.84461872214889
$ws.Range("R5").Value = 556.60156849934
$ws.Range("S5").Value = 0.05866910407468769
$ws.Range("T5").Value = 0.05866910407468769
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.77811266666667
$ws.Range("H6").Value = 32.334338
$ws.Range("I6").Value = 0.351787262211666
$ws.Range("J6").Value = 0.351787262211666
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.62053
$ws.Range("N6").Value = 22.86159
$ws.Range("O6").Value = 0.2214906134136664
$ws.Range("P6").Value = 0.2214906134136664
$ws.Range("Q6").Value = 82.13493091971333
$ws.Range("R6").Value = 739.21437827742
$ws.Range("S6").Value = 0.0779175764983762
$ws.Range("T6").Value = 0.07791757649837619
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.77811266666667
$ws.Range("H7").Value = 32.334338
$ws.Range("I7").Value = 0.351787262211666
$ws.Range("J7").Value = 0.351787262211666
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.340016666666665
$ws.Range("N7").Value = 28.02005
$ws.Range("O7").Value = 0.2714674728390108
$ws.Range("P7").Value = 0.2714674728390108
$ws.Range("Q7").Value = 100.6677519418778
$ws.Range("R7").Value = 906.0097674769
$ws.Range("S7").Value = 0.09549879904955542
$ws.Range("T7").Value = 0.09549879904955542
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.77811266666667
$ws.Range("H8").Value = 32.334338
$ws.Range("I8").Value = 0.351787262211666
$ws.Range("J8").Value = 0.351787262211666
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.970095666666666
$ws.Range("N8").Value = 23.910287
$ws.Range("O8").Value = 0.2316507353393536
$ws.Range("P8").Value = 0.2316507353393536
$ws.Range("Q8").Value = 85.90258905944511
$ws.Range("R8").Value = 773.123301535006
$ws.Range("S8").Value = 0.08149177797435044
$ws.Range("T8").Value = 0.08149177797435044
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.77811266666667
$ws.Range("H9").Value = 32.334338
$ws.Range("I9").Value = 0.351787262211666
$ws.Range("J9").Value = 0.351787262211666
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.475014333333332
$ws.Range("N9").Value = 28.425043
$ws.Range("O9").Value = 0.2753911784079691
$ws.Range("P9").Value = 0.2753911784079691
$ws.Range("Q9").Value = 102.1227720029482
$ws.Range("R9").Value = 919.104948026534
$ws.Range("S9").Value = 0.09687910868938393
$ws.Range("T9").Value = 0.09687910868938393
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.559571333333333
$ws.Range("H10").Value = 19.678714
$ws.Range("I10").Value = 0.2140981182885632
$ws.Range("J10").Value = 0.2140981182885631
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.62053
$ws.Range("N10").Value = 22.86159
$ws.Range("O10").Value = 0.2214906134136664
$ws.Range("P10").Value = 0.2214906134136664
$ws.Range("Q10").Value = 49.98741013280667
$ws.Range("R10").Value = 449.88669119526
$ws.Range("S10").Value = 0.04742072355044556
$ws.Range("T10").Value = 0.04742072355044555
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.559571333333333
$ws.Range("H11").Value = 19.678714
$ws.Range("I11").Value = 0.2140981182885632
$ws.Range("J11").Value = 0.2140981182885631
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.340016666666665
$ws.Range("N11").Value = 28.02005
$ws.Range("O11").Value = 0.2714674728390108
$ws.Range("P11").Value = 0.2714674728390108
$ws.Range("Q11").Value = 61.26650557952222
$ws.Range("R11").Value = 551.3985502157
$ws.Range("S11").Value = 0.05812067511138386
$ws.Range("T11").Value = 0.05812067511138385
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.559571333333333
$ws.Range("H12").Value = 19.678714
$ws.Range("I12").Value = 0.2140981182885632
$ws.Range("J12").Value = 0.2140981182885631
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.970095666666666
$ws.Range("N12").Value = 23.910287
$ws.Range("O12").Value = 0.2316507353393536
$ws.Range("P12").Value = 0.2316507353393536
$ws.Range("Q12").Value = 52.28041105899089
$ws.Range("R12").Value = 470.523699530918
$ws.Range("S12").Value = 0.04959598653631757
$ws.Range("T12").Value = 0.04959598653631756
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.559571333333333
$ws.Range("H13").Value = 19.678714
$ws.Range("I13").Value = 0.2140981182885632
$ws.Range("J13").Value = 0.2140981182885631
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.475014333333332
$ws.Range("N13").Value = 28.425043
$ws.Range("O13").Value = 0.2753911784079691
$ws.Range("P13").Value = 0.2753911784079691
$ws.Range("Q13").Value = 62.15203240385577
$ws.Range("R13").Value = 559.3682916347019
$ws.Range("S13").Value = 0.05896073309041618
$ws.Range("T13").Value = 0.05896073309041617
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.773344333333334
$ws.Range("H14").Value = 20.320033
$ws.Range("I14").Value = 0.2210754640197275
$ws.Range("J14").Value = 0.2210754640197275
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 7.62053
$ws.Range("N14").Value = 22.86159
$ws.Range("O14").Value = 0.2214906134136664
$ws.Range("P14").Value = 0.2214906134136664
$ws.Range("Q14").Value = 51.61647369249667
$ws.Range("R14").Value = 464.54826323247
$ws.Range("S14").Value = 0.04896614013644038
$ws.Range("T14").Value = 0.04896614013644036
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.773344333333334
$ws.Range("H15").Value = 20.320033
$ws.Range("I15").Value = 0.2210754640197275
$ws.Range("J15").Value = 0.2210754640197275
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.340016666666665
$ws.Range("N15").Value = 28.02005
$ws.Range("O15").Value = 0.2714674728390108
$ws.Range("P15").Value = 0.2714674728390108
$ws.Range("Q15").Value = 63.26314896240555
$ws.Range("R15").Value = 569.36834066165
$ws.Range("S15").Value = 0.0600147975241471
$ws.Range("T15").Value = 0.06001479752414709
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.773344333333334
$ws.Range("H16").Value = 20.320033
$ws.Range("I16").Value = 0.2210754640197275
$ws.Range("J16").Value = 0.2210754640197275
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.970095666666666
$ws.Range("N16").Value = 23.910287
$ws.Range("O16").Value = 0.2316507353393536
$ws.Range("P16").Value = 0.2316507353393536
$ws.Range("Q16").Value = 53.98420231994123
$ws.Range("R16").Value = 485.857820879471
$ws.Range("S16").Value = 0.05121229380565868
$ws.Range("T16").Value = 0.05121229380565868
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.773344333333334
$ws.Range("H17").Value = 20.320033
$ws.Range("I17").Value = 0.2210754640197275
$ws.Range("J17").Value = 0.2210754640197275
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.475014333333332
$ws.Range("N17").Value = 28.425043
$ws.Range("O17").Value = 0.2753911784079691
$ws.Range("P17").Value = 0.2753911784079691
$ws.Range("Q17").Value = 64.17753464293544
$ws.Range("R17").Value = 577.597811786419
$ws.Range("S17").Value = 0.06088223255348133
$ws.Range("T17").Value = 0.06088223255348132

Write-Output "Applied NATMI value updates"